$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.182.98"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").Value = "2.523.93"
$ws.Range("E3").Value = "  +0.91%  "

$ws.Range("E4").Value = "  +0.00%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "323.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "109.20"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("E7").Value = "  +0.80%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.558"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.06%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "40.59"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.67%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.37"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +10.75%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0823"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.23%  "

$ws.Range("E13").Value = "  +1.11%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.28"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").Value = "2.922.39"

$ws.Range("D16").Value = "2.524.71"
$ws.Range("E16").Value = "  +0.95%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.859"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").Value = "48.051.81"
$ws.Range("E18").Value = "  +1.72%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "13.28"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +3.28%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.64"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("E22").Value = "  +0.33%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "72.44"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.38%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "269.90"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +8.91%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.21"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("E27").Value = "  +0.03%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.16"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.88%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.146"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +5.56%  "

$ws.Range("E30").Value = "  -3.89%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "35.75"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.42%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "49.84"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.02%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.97"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.41"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  -0.04%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0795"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.77%  "

$ws.Range("E37").Value = "  +0.88%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.75"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.03%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.99"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("E40").Value = "  -0.08%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "22.32"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +4.87%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.00%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "119.32"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").Value = "2.013.84"
$ws.Range("E45").Value = "  +1.12%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.14"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.47%  "

$ws.Range("E47").Value = "  -0.67%  "

$ws.Range("E48").Value = "  +5.45%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.14"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "80.08"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.56%  "
